# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for the affected leve rows on each job sheet.
# Values come from the upstream market API snapshot; only H:N change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 589.4
$ws.Range("I41").Value = 366.8
$ws.Range("J41").Value = 812
$ws.Range("K41").Value = 366.8
$ws.Range("L41").Value = 812
$ws.Range("M41").Value = 73.19999999999999
$ws.Range("N41").Value = -1692

$ws.Range("H58").Value = 1017.3571
$ws.Range("I58").Value = 353.58334
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 1060.75002
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -910.7500199999999
$ws.Range("N58").Value = -15300

$ws.Range("H62").Value = 3124.5
$ws.Range("I62").Value = 2250
$ws.Range("K62").Value = 2250
$ws.Range("M62").Value = -1626

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H65").Value = 3124.5
$ws.Range("I65").Value = 2250
$ws.Range("K65").Value = 11250
$ws.Range("M65").Value = -8130

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H74").Value = 3899
$ws.Range("I74").Value = 3499
$ws.Range("K74").Value = 3499
$ws.Range("M74").Value = -2563

$ws.Range("H77").Value = 3899
$ws.Range("I77").Value = 3499
$ws.Range("K77").Value = 17495
$ws.Range("M77").Value = -12815

$ws.Range("H121").Value = 998
$ws.Range("J121").Value = 998
$ws.Range("L121").Value = 2994
$ws.Range("N121").Value = -6488

$ws.Range("H132").Value = 1104.6072
$ws.Range("I132").Value = 1104.6072
$ws.Range("K132").Value = 3313.8216
$ws.Range("M132").Value = -783.8215999999998

$ws.Range("H137").Value = 1998.6666
$ws.Range("I137").Value = 1300
$ws.Range("K137").Value = 3900
$ws.Range("M137").Value = -1350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5099.3413
$ws.Range("I32").Value = 3585.7144
$ws.Range("J32").Value = 8359.462
$ws.Range("K32").Value = 3585.7144
$ws.Range("L32").Value = 8359.462
$ws.Range("M32").Value = -3298.7144
$ws.Range("N32").Value = -8933.462

$ws.Range("H62").Value = 29999
$ws.Range("J62").Value = 29999
$ws.Range("L62").Value = 29999
$ws.Range("N62").Value = -31247

$ws.Range("H63").Value = 1866
$ws.Range("I63").Value = 1866
$ws.Range("K63").Value = 1866
$ws.Range("M63").Value = -1180

$ws.Range("H65").Value = 29999
$ws.Range("J65").Value = 29999
$ws.Range("L65").Value = 89997
$ws.Range("N65").Value = -96237

$ws.Range("H66").Value = 1866
$ws.Range("I66").Value = 1866
$ws.Range("K66").Value = 9330
$ws.Range("M66").Value = -5898

$ws.Range("H122").Value = 1659.6364
$ws.Range("I122").Value = 1477.1538
$ws.Range("K122").Value = 4431.4614
$ws.Range("M122").Value = -1981.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1800.0714
$ws.Range("I20").Value = 1461.2
$ws.Range("J20").Value = 2647.25
$ws.Range("K20").Value = 1461.2
$ws.Range("L20").Value = 2647.25
$ws.Range("M20").Value = -1214.2
$ws.Range("N20").Value = -3141.25

$ws.Range("H76").Value = 45077.75
$ws.Range("J76").Value = 45077.75
$ws.Range("L76").Value = 45077.75
$ws.Range("N76").Value = -45707.75

$ws.Range("H79").Value = 45077.75
$ws.Range("J79").Value = 45077.75
$ws.Range("L79").Value = 45077.75
$ws.Range("N79").Value = -47261.75

$ws.Range("H95").Value = 71898
$ws.Range("J95").Value = 71898
$ws.Range("L95").Value = 71898
$ws.Range("N95").Value = -77390

$ws.Range("H99").Value = 899
$ws.Range("I99").Value = 899
$ws.Range("K99").Value = 899
$ws.Range("M99").Value = 599

$ws.Range("H100").Value = 22250
$ws.Range("J100").Value = 22250
$ws.Range("L100").Value = 22250
$ws.Range("N100").Value = -24414

$ws.Range("H107").Value = 2134.4285
$ws.Range("I107").Value = 1990.1666
$ws.Range("K107").Value = 1990.1666
$ws.Range("M107").Value = -70.16660000000002

$ws.Range("H134").Value = 7091.023
$ws.Range("I134").Value = 7624.212
$ws.Range("J134").Value = 5491.4546
$ws.Range("K134").Value = 22872.636
$ws.Range("L134").Value = 16474.3638
$ws.Range("M134").Value = -20337.636
$ws.Range("N134").Value = -21544.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747

$ws.Range("H134").Value = 905.5
$ws.Range("I134").Value = 749.1429000000001
$ws.Range("K134").Value = 2247.4287
$ws.Range("M134").Value = 287.5712999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 81.2
$ws.Range("I33").Value = 97.71429000000001
$ws.Range("K33").Value = 586.28574
$ws.Range("M33").Value = -303.28574

$ws.Range("H131").Value = 762
$ws.Range("J131").Value = 804.4945
$ws.Range("L131").Value = 2413.4835
$ws.Range("N131").Value = -12493.4835

$ws.Range("H132").Value = 1475.7693
$ws.Range("I132").Value = 1462.2727
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 13160.4543
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -10630.4543
$ws.Range("N132").Value = -19010

$ws.Range("H138").Value = 2280.8333
$ws.Range("I138").Value = 1811.6666
$ws.Range("J138").Value = 2750
$ws.Range("K138").Value = 5434.9998
$ws.Range("L138").Value = 8250
$ws.Range("M138").Value = -294.9997999999996
$ws.Range("N138").Value = -18530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2299.9
$ws.Range("I97").Value = 1983.3334
$ws.Range("K97").Value = 1983.3334
$ws.Range("M97").Value = -1487.3334

$ws.Range("H122").Value = 2366.8572
$ws.Range("I122").Value = 2254.2727
$ws.Range("K122").Value = 6762.8181
$ws.Range("M122").Value = -4312.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3750
$ws.Range("I61").Value = 3333.3333
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3333.3333
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3131.3333
$ws.Range("N61").Value = -5404

$ws.Range("H93").Value = 538.3333
$ws.Range("I93").Value = 356.375
$ws.Range("K93").Value = 356.375
$ws.Range("M93").Value = 891.625

$ws.Range("H113").Value = 3750
$ws.Range("I113").Value = 3333.3333
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3333.3333
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1163.3333
$ws.Range("N113").Value = -9340

$ws.Range("H132").Value = 2282.4348
$ws.Range("I132").Value = 1719.5
$ws.Range("K132").Value = 5158.5
$ws.Range("M132").Value = -2628.5

$ws.Range("H136").Value = 2993.1482
$ws.Range("I136").Value = 1735.9375
$ws.Range("K136").Value = 5207.8125
$ws.Range("M136").Value = -2657.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 28592
$ws.Range("J92").Value = 28592
$ws.Range("L92").Value = 28592
$ws.Range("N92").Value = -33584

$ws.Range("H100").Value = 199.5
$ws.Range("I100").Value = 199.5
$ws.Range("K100").Value = 399
$ws.Range("M100").Value = 142

$ws.Range("H113").Value = 647
$ws.Range("I113").Value = 522.7857
$ws.Range("J113").Value = 1226.6666
$ws.Range("K113").Value = 1568.3571
$ws.Range("L113").Value = 3679.9998
$ws.Range("M113").Value = 601.6428999999998
$ws.Range("N113").Value = -8019.9998

$ws.Range("H122").Value = 98957.625
$ws.Range("I122").Value = 157132.2
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 471396.6
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -468946.6
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 8439.76
$ws.Range("I132").Value = 3999
$ws.Range("K132").Value = 11997
$ws.Range("M132").Value = -9467

$ws.Range("H136").Value = 2515.3635
$ws.Range("I136").Value = 1368
$ws.Range("K136").Value = 4104
$ws.Range("M136").Value = -1554
